# "Adding The Day Work"
# Appends new day-tracking entries (rows 26-43) to the work log, mirroring
# the existing repeating block pattern (a merged date cell in column A,
# spanning 6 rows, with TASK / time / status detail rows beneath it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Block 1: 2020-11-27 (rows 26:31)
# ---------------------------------------------------------------------
# Merge first (merging afterwards would stamp every cell in the range
# with the top-left cell's style, wiping out the per-row formatting
# pattern we are about to paste in).
$ws.Range("A26:A31").Merge()

$ws.Range("A26").Formula = "=DATE(2020,11,27)"
$ws.Range("B26").Value = "Radiante ApI working"
$ws.Range("C26").Value = "9:00:00 PM : 11:00 PM"
$ws.Range("D26").Value = "in progress"

# Match the existing look of the other date blocks: copy the A-column
# format pattern (merged "date" cell + blank filler cells) down from the
# 2020-10-17 block, and copy the time-value format onto the new C/D cells.
$ws.Range("A20:A25").Copy()
$ws.Range("A26:A31").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

$ws.Range("C14").Copy()
$ws.Range("C26:D26").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# Block 2: 2020-11-29 (rows 32:37)
# ---------------------------------------------------------------------
$ws.Range("A32:A37").Merge()

$ws.Range("A32").Formula = "=DATE(2020,11,29)"
$ws.Range("B32").Value = "Quraan 2 Quarters"
$ws.Range("C32").Value = "5:00/5:30"

$ws.Range("B33").Value = "Fajr"
$ws.Range("C33").Value = "4:30/5:00"

$ws.Range("B34").Value = "Radiante ApI working"
$ws.Range("C34").Value = "5:30/7:00"

$ws.Range("B35").Value = "Travelling"
$ws.Range("C35").Value = "7:00/9:00"

$ws.Range("B36").Value = "Radiante ApI working"
$ws.Range("C36").Value = "9:00/11:40"

$ws.Range("B37").Value = "Duhr"

$ws.Range("A20:A25").Copy()
$ws.Range("A32:A37").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# Block 3: trailing date-only block (rows 38:43)
# ---------------------------------------------------------------------
$ws.Range("A38:A43").Merge()

$ws.Range("A38").Formula = "=DATE(2020,11,29)"

$ws.Range("A20:A25").Copy()
$ws.Range("A38:A43").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# View state: scroll/selection left where the author last worked.
# ---------------------------------------------------------------------
$ws.Range("A16").Select()
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("C37").Select()
